$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3069.9062
$ws.Range("I64").Value = 2841.0527
$ws.Range("J64").Value = 3404.3845
$ws.Range("K64").Value = 2841.0527
$ws.Range("L64").Value = 3404.3845
$ws.Range("M64").Value = -2593.0527
$ws.Range("N64").Value = -3900.3845
$ws.Range("H67").Value = 3069.9062
$ws.Range("I67").Value = 2841.0527
$ws.Range("J67").Value = 3404.3845
$ws.Range("K67").Value = 2841.0527
$ws.Range("L67").Value = 3404.3845
$ws.Range("M67").Value = -1983.0527
$ws.Range("N67").Value = -5120.3845
$ws.Range("H70").Value = 2601.1924
$ws.Range("I70").Value = 951
$ws.Range("J70").Value = 2901.2273
$ws.Range("K70").Value = 2853
$ws.Range("L70").Value = 8703.6819
$ws.Range("M70").Value = -2583
$ws.Range("N70").Value = -9243.6819
$ws.Range("H73").Value = 2601.1924
$ws.Range("I73").Value = 951
$ws.Range("J73").Value = 2901.2273
$ws.Range("K73").Value = 2853
$ws.Range("L73").Value = 8703.6819
$ws.Range("M73").Value = -1917
$ws.Range("N73").Value = -10575.6819
$ws.Range("H74").Value = 4067.5
$ws.Range("I74").Value = 2990.7693
$ws.Range("J74").Value = 8733.333000000001
$ws.Range("K74").Value = 2990.7693
$ws.Range("L74").Value = 8733.333000000001
$ws.Range("M74").Value = -2054.7693
$ws.Range("N74").Value = -10605.333
$ws.Range("H77").Value = 4067.5
$ws.Range("I77").Value = 2990.7693
$ws.Range("J77").Value = 8733.333000000001
$ws.Range("K77").Value = 14953.8465
$ws.Range("L77").Value = 43666.665
$ws.Range("M77").Value = -10273.8465
$ws.Range("N77").Value = -53026.665
$ws.Range("H141").Value = 2460.484
$ws.Range("I141").Value = 687.12
$ws.Range("K141").Value = 2061.36
$ws.Range("M141").Value = 3118.64

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 9
$ws.Range("I29").Value = 9
$ws.Range("K29").Value = 9
$ws.Range("M29").Value = 299
$ws.Range("H32").Value = 7769
$ws.Range("I32").Value = 6455.4346
$ws.Range("J32").Value = 22875
$ws.Range("K32").Value = 6455.4346
$ws.Range("L32").Value = 22875
$ws.Range("M32").Value = -6168.4346
$ws.Range("N32").Value = -23449

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 16755885
$ws.Range("I134").Value = 24709828
$ws.Range("J134").Value = 52606.4
$ws.Range("K134").Value = 74129484
$ws.Range("L134").Value = 157819.2
$ws.Range("M134").Value = -74126949
$ws.Range("N134").Value = -162889.2

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 681.2692
$ws.Range("I16").Value = 675.4211
$ws.Range("J16").Value = 697.1429000000001
$ws.Range("K16").Value = 675.4211
$ws.Range("L16").Value = 697.1429000000001
$ws.Range("M16").Value = -388.4211
$ws.Range("N16").Value = -1271.1429
$ws.Range("H31").Value = 4171803.2
$ws.Range("I31").Value = 7577198.5
$ws.Range("J31").Value = 9653.333000000001
$ws.Range("K31").Value = 7577198.5
$ws.Range("L31").Value = 9653.333000000001
$ws.Range("M31").Value = -7576903.5
$ws.Range("N31").Value = -10243.333
$ws.Range("H34").Value = 4171803.2
$ws.Range("I34").Value = 7577198.5
$ws.Range("J34").Value = 9653.333000000001
$ws.Range("K34").Value = 7577198.5
$ws.Range("L34").Value = 9653.333000000001
$ws.Range("M34").Value = -7576996.5
$ws.Range("N34").Value = -10057.333
$ws.Range("H62").Value = 2278.4285
$ws.Range("I62").Value = 2235.125
$ws.Range("J62").Value = 2336.1667
$ws.Range("K62").Value = 2235.125
$ws.Range("L62").Value = 2336.1667
$ws.Range("M62").Value = -1611.125
$ws.Range("N62").Value = -3584.1667
$ws.Range("H65").Value = 2278.4285
$ws.Range("I65").Value = 2235.125
$ws.Range("J65").Value = 2336.1667
$ws.Range("K65").Value = 11175.625
$ws.Range("L65").Value = 11680.8335
$ws.Range("M65").Value = -8055.625
$ws.Range("N65").Value = -17920.8335
$ws.Range("H113").Value = 681.2692
$ws.Range("I113").Value = 675.4211
$ws.Range("J113").Value = 697.1429000000001
$ws.Range("K113").Value = 675.4211
$ws.Range("L113").Value = 697.1429000000001
$ws.Range("M113").Value = 1494.5789
$ws.Range("N113").Value = -5037.1429

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2495.1924
$ws.Range("J39").Value = 2738.0435
$ws.Range("L39").Value = 8214.130500000001
$ws.Range("N39").Value = -8802.130500000001
$ws.Range("H105").Value = 7282.6
$ws.Range("J105").Value = 7422.222
$ws.Range("L105").Value = 22266.666
$ws.Range("N105").Value = -27508.666
$ws.Range("H129").Value = 2038.9032
$ws.Range("I129").Value = 726.6667
$ws.Range("J129").Value = 3855.8462
$ws.Range("K129").Value = 2180.0001
$ws.Range("L129").Value = 11567.5386
$ws.Range("M129").Value = 2819.9999
$ws.Range("N129").Value = -21567.5386
$ws.Range("H140").Value = 2398.946
$ws.Range("I140").Value = 2133.1614
$ws.Range("J140").Value = 3772.1667
$ws.Range("K140").Value = 6399.4842
$ws.Range("L140").Value = 11316.5001
$ws.Range("M140").Value = -1219.4842
$ws.Range("N140").Value = -21676.5001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 17666.666
$ws.Range("J40").Value = 17666.666
$ws.Range("L40").Value = 17666.666
$ws.Range("N40").Value = -17968.666
$ws.Range("H47").Value = 12510.5
$ws.Range("J47").Value = 12510.5
$ws.Range("L47").Value = 12510.5
$ws.Range("N47").Value = -13646.5
$ws.Range("H70").Value = 6227132.5
$ws.Range("I70").Value = 2608516.8
$ws.Range("J70").Value = 15876774
$ws.Range("K70").Value = 2608516.8
$ws.Range("L70").Value = 15876774
$ws.Range("M70").Value = -2608246.8
$ws.Range("N70").Value = -15877314
$ws.Range("H73").Value = 6227132.5
$ws.Range("I73").Value = 2608516.8
$ws.Range("J73").Value = 15876774
$ws.Range("K73").Value = 2608516.8
$ws.Range("L73").Value = 15876774
$ws.Range("M73").Value = -2607580.8
$ws.Range("N73").Value = -15878646
$ws.Range("H80").Value = 13500
$ws.Range("I80").Value = 5581.8184
$ws.Range("J80").Value = 30920
$ws.Range("K80").Value = 5581.8184
$ws.Range("L80").Value = 30920
$ws.Range("M80").Value = -4583.8184
$ws.Range("N80").Value = -32916
$ws.Range("H83").Value = 13500
$ws.Range("I83").Value = 5581.8184
$ws.Range("J83").Value = 30920
$ws.Range("K83").Value = 27909.092
$ws.Range("L83").Value = 154600
$ws.Range("M83").Value = -22917.092
$ws.Range("N83").Value = -164584

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8335500
$ws.Range("I136").Value = 31250752
$ws.Range("J136").Value = 2680.4546
$ws.Range("K136").Value = 93752256
$ws.Range("L136").Value = 8041.3638
$ws.Range("M136").Value = -93749706
$ws.Range("N136").Value = -13141.3638
